$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 423
$ws.Range("E2").Value = -5
$ws.Range("F2").Value = -5
$ws.Range("G2").Value = 8
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1453
$ws.Range("L2").Value = 302
$ws.Range("M2").Value = 1151
$ws.Range("N2").Value = 1135
$ws.Range("O2").Value = 16
$ws.Range("P2").Value = 123
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 43
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 8
$ws.Range("W2").Value = -1.23
$ws.Range("X2").Value = 1.98
$ws.Range("Y2").Value = 0.6899999999999999
$ws.Range("Z2").Value = 0.59
$ws.Range("AA2").Value = 26.24
$ws.Range("AB2").Value = 821.41
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 443.27
$ws.Range("AE2").Value = 2778
$ws.Range("AF2").Value = 3.03
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 40836231

# Row 3
$ws.Range("D3").Value = 447
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = -500
$ws.Range("H3").Value = -499
$ws.Range("I3").Value = -500
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2490
$ws.Range("L3").Value = 838
$ws.Range("M3").Value = 1652
$ws.Range("N3").Value = 1635
$ws.Range("O3").Value = 17
$ws.Range("P3").Value = 151
$ws.Range("Q3").Value = -2
$ws.Range("R3").Value = -20
$ws.Range("S3").Value = 9
$ws.Range("T3").Value = 15
$ws.Range("U3").Value = -17
$ws.Range("V3").Value = 17
$ws.Range("W3").Value = 3.44
$ws.Range("X3").Value = -111.7
$ws.Range("Y3").Value = -36.07
$ws.Range("Z3").Value = -25.31
$ws.Range("AA3").Value = 50.7
$ws.Range("AB3").Value = 980.77
$ws.Range("AC3").Value = -1171
$ws.Range("AD3").Value = -6.57
$ws.Range("AE3").Value = 3257
$ws.Range("AF3").Value = 2.36
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 50209620

# Row 4
$ws.Range("D4").Value = 504
$ws.Range("E4").Value = 33
$ws.Range("F4").Value = 33
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 15
$ws.Range("I4").Value = 15
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3622
$ws.Range("L4").Value = 1961
$ws.Range("M4").Value = 1661
$ws.Range("N4").Value = 1645
$ws.Range("O4").Value = 17
$ws.Range("P4").Value = 227
$ws.Range("Q4").Value = 31
$ws.Range("R4").Value = -87
$ws.Range("S4").Value = 109
$ws.Range("T4").Value = 89
$ws.Range("U4").Value = -58
$ws.Range("V4").Value = 127
$ws.Range("W4").Value = 6.54
$ws.Range("X4").Value = 3.02
$ws.Range("Y4").Value = 0.93
$ws.Range("Z4").Value = 0.5
$ws.Range("AA4").Value = 118.03
$ws.Range("AB4").Value = 624.58
$ws.Range("AC4").Value = 31
$ws.Range("AD4").Value = 218.25
$ws.Range("AE4").Value = 3276
$ws.Range("AF4").Value = 2.03
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 50209620

# Row 5
$ws.Range("D5").Value = 692
$ws.Range("E5").Value = 51
$ws.Range("F5").Value = 51
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 34
$ws.Range("I5").Value = 33
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4018
$ws.Range("L5").Value = 2293
$ws.Range("M5").Value = 1725
$ws.Range("N5").Value = 1707
$ws.Range("O5").Value = 17
$ws.Range("P5").Value = 227
$ws.Range("Q5").Value = 65
$ws.Range("R5").Value = -32
$ws.Range("S5").Value = 299
$ws.Range("T5").Value = 51
$ws.Range("U5").Value = 14
$ws.Range("V5").Value = 404
$ws.Range("W5").Value = 7.31
$ws.Range("X5").Value = 4.88
$ws.Range("Y5").Value = 1.99
$ws.Range("Z5").Value = 0.89
$ws.Range("AA5").Value = 132.97
$ws.Range("AB5").Value = 652.28
$ws.Range("AC5").Value = 67
$ws.Range("AD5").Value = 146.77
$ws.Range("AE5").Value = 3401
$ws.Range("AF5").Value = 2.87
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 50209620

# Row 6
$ws.Range("D6").Value = 761
$ws.Range("E6").Value = 30
$ws.Range("F6").Value = 30
$ws.Range("G6").Value = -1017
$ws.Range("H6").Value = -1129
$ws.Range("I6").Value = -1085
$ws.Range("K6").Value = 5068
$ws.Range("L6").Value = 2258
$ws.Range("M6").Value = 2810
$ws.Range("N6").Value = 2836
$ws.Range("P6").Value = 326
$ws.Range("Q6").Value = 21
$ws.Range("R6").Value = -2119
$ws.Range("S6").Value = 1765
$ws.Range("T6").Value = 98
$ws.Range("U6").Value = -78
$ws.Range("V6").Value = 294
$ws.Range("W6").Value = 3.95
$ws.Range("X6").Value = -148.37
$ws.Range("Y6").Value = -47.74
$ws.Range("Z6").Value = -24.85
$ws.Range("AA6").Value = 80.38
$ws.Range("AB6").Value = 757.26
$ws.Range("AC6").Value = -2048
$ws.Range("AD6").Value = -7.1
$ws.Range("AE6").Value = 4353
$ws.Range("AF6").Value = 3.34
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 65163174

# Clear D7:AI9 (rows 7-9 retain only A,B,C)
$ws.Range("D7:AI9").ClearContents()

# AI6 is removed entirely (no longer present)
$ws.Range("AI6").ClearContents()
